# Auto-generated script applying cached-value updates described in the commit diff.
# The workbook stores plain (non-formula) cached numeric values in columns H-N of
# each profession sheet, so the edit is simply a set of direct cell value writes
# (plus two cells whose content is removed entirely).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 83190
$ws.Range("J17").Value = 91389
$ws.Range("L17").Value = 274167
$ws.Range("N17").Value = -274503
$ws.Range("H38").Value = 474.84616
$ws.Range("I38").Value = 347.75
$ws.Range("J38").Value = 2000
$ws.Range("K38").Value = 1043.25
$ws.Range("L38").Value = 6000
$ws.Range("M38").Value = -671.25
$ws.Range("N38").Value = -6744
$ws.Range("H39").Value = 1337.75
$ws.Range("I39").Value = 1393.6
$ws.Range("K39").Value = 4180.799999999999
$ws.Range("M39").Value = -3884.799999999999
$ws.Range("H42").Value = 269.1
$ws.Range("I42").Value = 59
$ws.Range("J42").Value = 479.2
$ws.Range("K42").Value = 177
$ws.Range("L42").Value = 1437.6
$ws.Range("M42").Value = 53
$ws.Range("N42").Value = -1897.6
$ws.Range("H43").Value = 9912.25
$ws.Range("I43").Value = 9499.666999999999
$ws.Range("K43").Value = 9499.666999999999
$ws.Range("M43").Value = -9430.666999999999
$ws.Range("H49").Value = 1321.8
$ws.Range("I49").Value = 165
$ws.Range("J49").Value = 2093
$ws.Range("K49").Value = 495
$ws.Range("L49").Value = 6279
$ws.Range("M49").Value = -359
$ws.Range("N49").Value = -6551
$ws.Range("H74").Value = 4022.9285
$ws.Range("I74").Value = 3040.125
$ws.Range("K74").Value = 3040.125
$ws.Range("M74").Value = -2104.125
$ws.Range("H77").Value = 4022.9285
$ws.Range("I77").Value = 3040.125
$ws.Range("K77").Value = 15200.625
$ws.Range("M77").Value = -10520.625
$ws.Range("H86").Value = 2305.75
$ws.Range("I86").Value = 2149.5
$ws.Range("K86").Value = 2149.5
$ws.Range("M86").Value = -1026.5
$ws.Range("H89").Value = 2305.75
$ws.Range("I89").Value = 2149.5
$ws.Range("K89").Value = 10747.5
$ws.Range("M89").Value = -5131.5
$ws.Range("H106").Value = 12196.9
$ws.Range("I106").Value = 1995.4
$ws.Range("K106").Value = 1995.4
$ws.Range("M106").Value = -1364.4
$ws.Range("H137").Value = 4474.647
$ws.Range("I137").Value = 1556.5172
$ws.Range("K137").Value = 4669.5516
$ws.Range("M137").Value = -2119.5516
$ws.Range("H138").Value = 2307.5386
$ws.Range("I138").Value = 1421.1052
$ws.Range("J138").Value = 3149.65
$ws.Range("K138").Value = 4263.3156
$ws.Range("L138").Value = 9448.950000000001
$ws.Range("M138").Value = 876.6844000000001
$ws.Range("N138").Value = -19728.95
$ws.Range("H141").Value = 45674.094
$ws.Range("I141").Value = 59456.875
$ws.Range("J141").Value = 1569.2
$ws.Range("K141").Value = 178370.625
$ws.Range("L141").Value = 4707.6
$ws.Range("M141").Value = -173190.625
$ws.Range("N141").Value = -15067.6

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 42533.98
$ws.Range("I32").Value = 23627.645
$ws.Range("J32").Value = 184331.5
$ws.Range("K32").Value = 23627.645
$ws.Range("L32").Value = 184331.5
$ws.Range("M32").Value = -23340.645
$ws.Range("N32").Value = -184905.5
$ws.Range("H45").Value = 597287.7
$ws.Range("I45").Value = 1265835.1
$ws.Range("J45").Value = 3023.3333
$ws.Range("K45").Value = 1265835.1
$ws.Range("L45").Value = 3023.3333
$ws.Range("M45").Value = -1265458.1
$ws.Range("N45").Value = -3777.3333
$ws.Range("H74").Value = 1408.3077
$ws.Range("I74").Value = 1362.4286
$ws.Range("J74").Value = 1601
$ws.Range("K74").Value = 1362.4286
$ws.Range("L74").Value = 1601
$ws.Range("M74").Value = -488.4286
$ws.Range("N74").Value = -3349
$ws.Range("H77").Value = 1408.3077
$ws.Range("I77").Value = 1362.4286
$ws.Range("J77").Value = 1601
$ws.Range("K77").Value = 6812.143
$ws.Range("L77").Value = 8005
$ws.Range("M77").Value = -2444.143
$ws.Range("N77").Value = -16741

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 9065.096
$ws.Range("I20").Value = 7762.7144
$ws.Range("K20").Value = 7762.7144
$ws.Range("M20").Value = -7515.7144

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3297.6924
$ws.Range("I132").Value = 2876.7778
$ws.Range("K132").Value = 8630.3334
$ws.Range("M132").Value = -6100.3334

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 125349.875
$ws.Range("I50").Value = 250149.75
$ws.Range("K50").Value = 750449.25
$ws.Range("M50").Value = -749968.25
$ws.Range("H53").Value = 125349.875
$ws.Range("I53").Value = 250149.75
$ws.Range("K53").Value = 750449.25
$ws.Range("M53").Value = -749968.25
$ws.Range("H128").Value = 373905.75
$ws.Range("I128").Value = 373905.75
$ws.Range("K128").Value = 1121717.25
$ws.Range("M128").Value = -1116737.25

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 25999.25
$ws.Range("I43").Value = 1999.5
$ws.Range("K43").Value = 1999.5
$ws.Range("M43").Value = -1848.5
$ws.Range("H102").Value = 2806.1052
$ws.Range("I102").Value = 2518.875
$ws.Range("K102").Value = 2518.875
$ws.Range("M102").Value = -896.875
$ws.Range("H127").Value = 66666.336
$ws.Range("I127").Value = 60000
$ws.Range("J127").Value = 69999.5
$ws.Range("K127").Value = 60000
$ws.Range("L127").Value = 69999.5
$ws.Range("M127").Value = -55040
$ws.Range("N127").Value = -79919.5
$ws.Range("H132").Value = 2353.9033
$ws.Range("J132").Value = 3259
$ws.Range("L132").Value = 9777
$ws.Range("N132").Value = -14837

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 32083.334
$ws.Range("I7").Value = 53370
$ws.Range("J7").Value = 5475
$ws.Range("K7").Value = 53370
$ws.Range("L7").Value = 5475
$ws.Range("M7").Value = -53258
$ws.Range("N7").Value = -5699
$ws.Range("H46").Value = 2146
$ws.Range("I46").Value = 2146
$ws.Range("K46").Value = 2146
$ws.Range("M46").Value = -1958
$ws.Range("H100").Value = 5574.1904
$ws.Range("I100").Value = 6080.643
$ws.Range("K100").Value = 6080.643
$ws.Range("M100").Value = -5539.643
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("H126").Value = 32083.334
$ws.Range("I126").Value = 53370
$ws.Range("J126").Value = 5475
$ws.Range("K126").Value = 160110
$ws.Range("L126").Value = 16425
$ws.Range("M126").Value = -157640
$ws.Range("N126").Value = -21365
$ws.Range("N110").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("H132").Value = 54569.215
$ws.Range("I132").Value = 55831.168
$ws.Range("J132").Value = 46997.5
$ws.Range("K132").Value = 167493.504
$ws.Range("L132").Value = 140992.5
$ws.Range("M132").Value = -164963.504
$ws.Range("N132").Value = -146052.5
$ws.Range("N112").ClearContents()
